# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
# Swap the match-record content (columns B..AD) between the given row pairs,
# leaving the row index (A), Div (C) and Date (D) columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(221, 222),
    @(282, 283),
    @(313, 314)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    # NOTE: reading the ".Value" property on this COM shim does not return
    # the actual cell contents (it yields the property descriptor instead),
    # so ".Value2" is used for both reading and writing here - it behaves
    # like the classic Value2 COM property and round-trips arrays fine.
    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
